# Financials update: insert a new "most recent period" column before column D
# (shifting the existing D:K data to E:L) and populate the new column D with
# the latest reporting period's figures, across all three statements
# (Income Statement, Balance Sheet, Cash Flow Statement).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column at D; Excel shifts D:K -> E:L automatically.
$ws.Columns("D").Insert()

# 2) The newly inserted column doesn't inherit the neighbouring column's
#    number formatting automatically here, so copy formats from column E
#    (which now holds what used to be column D) across the used rows.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)

# 3) Populate the new column D with the latest period's values, row by row.

# -- Income Statement --
$ws.Range("D7").Value2 = 43465
$ws.Range("D8").Value2 = 818300
$ws.Range("D9").Value2 = 695900
$ws.Range("D10").Value2 = 122400
$ws.Range("D12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = 200
$ws.Range("D15").Value2 = 0
$ws.Range("D17").Value2 = 759400
$ws.Range("D18").Value2 = 59000
$ws.Range("D20").Value2 = 16100
$ws.Range("D21").Value2 = 82100
$ws.Range("D22").Value2 = 7400
$ws.Range("D23").Value2 = 67600
$ws.Range("D24").Value2 = 5700
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = 61900
$ws.Range("D27").Value2 = 61900
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = 2300
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = -16100
$ws.Range("D33").Value2 = 64200
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = 64200

# -- Balance Sheet --
$ws.Range("D38").Value2 = 43465
$ws.Range("D41").Value2 = 30100
$ws.Range("D42").Value2 = 0
$ws.Range("D43").Value2 = 152600
$ws.Range("D44").Value2 = 0
$ws.Range("D45").Value2 = 0
$ws.Range("D46").Value2 = 0
$ws.Range("D47").Value2 = 2259100
$ws.Range("D48").Value2 = 54700
$ws.Range("D49").Value2 = 54500
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 0
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 3105100
$ws.Range("D57").Value2 = 0
$ws.Range("D58").Value2 = 0
$ws.Range("D59").Value2 = 1601400
$ws.Range("D60").Value2 = 0
$ws.Range("D61").Value2 = 149100
$ws.Range("D62").Value2 = 24200
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 2298200
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = 908100
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = 806800
$ws.Range("D77").Value2 = 0

# -- Cash Flow Statement --
$ws.Range("D80").Value2 = 43465
$ws.Range("D81").Value2 = 64200
$ws.Range("D83").Value2 = 7000
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = 217100
$ws.Range("D91").Value2 = -6100
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = -134200
$ws.Range("D96").Value2 = -83100
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = -77000
$ws.Range("D101").Value2 = 0
$ws.Range("D102").Value2 = 5900
